# Titan_Profits.xlsx refresh: scheduled-runner price/profit update
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H-N) for a batch of leve rows across all item-category sheets.
$wb = $excel.ActiveWorkbook

# ALC!row17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1389895.5
$ws.Range("J17").Value = 1538652.1
$ws.Range("L17").Value = 4615956.300000001
$ws.Range("N17").Value = -4616292.300000001

# ALC!row124
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H124").Value = 35000
$ws.Range("J124").Value = 35000
$ws.Range("L124").Value = 35000
$ws.Range("N124").Value = -44820

# ALC!row135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1312.9667
$ws.Range("I135").Value = 1245.375
$ws.Range("K135").Value = 11208.375
$ws.Range("M135").Value = -8673.375

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 8236681
$ws.Range("I138").Value = 2841935.8
$ws.Range("J138").Value = 11113879
$ws.Range("K138").Value = 8525807.399999999
$ws.Range("L138").Value = 33341637
$ws.Range("M138").Value = -8520667.399999999
$ws.Range("N138").Value = -33351917

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1237.1
$ws.Range("I45").Value = 1207.8889
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 1207.8889
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -830.8888999999999
$ws.Range("N45").Value = -2254

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2449.1282
$ws.Range("I61").Value = 1781.9333
$ws.Range("J61").Value = 4673.1113
$ws.Range("K61").Value = 1781.9333
$ws.Range("L61").Value = 4673.1113
$ws.Range("M61").Value = -1569.9333
$ws.Range("N61").Value = -5097.1113

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2005.2
$ws.Range("I122").Value = 1984.8889
$ws.Range("J122").Value = 2035.6666
$ws.Range("K122").Value = 5954.6667
$ws.Range("L122").Value = 6106.9998
$ws.Range("M122").Value = -3504.6667
$ws.Range("N122").Value = -11006.9998

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2449.1282
$ws.Range("I136").Value = 1781.9333
$ws.Range("J136").Value = 4673.1113
$ws.Range("K136").Value = 5345.7999
$ws.Range("L136").Value = 14019.3339
$ws.Range("M136").Value = -2795.7999
$ws.Range("N136").Value = -19119.3339

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 297274
$ws.Range("I105").Value = 3102.3809
$ws.Range("J105").Value = 772474.3
$ws.Range("K105").Value = 3102.3809
$ws.Range("L105").Value = 772474.3
$ws.Range("M105").Value = -1355.3809
$ws.Range("N105").Value = -775968.3

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 25643714
$ws.Range("I134").Value = 43480140
$ws.Range("J134").Value = 3851.5
$ws.Range("K134").Value = 130440420
$ws.Range("L134").Value = 11554.5
$ws.Range("M134").Value = -130437885
$ws.Range("N134").Value = -16624.5

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12953.182
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 12953.182
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 12953.182
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -13543.182

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 12953.182
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 12953.182
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 12953.182
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -13357.182

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2981486
$ws.Range("I99").Value = 4471367.5
$ws.Range("J99").Value = 1723.4286
$ws.Range("K99").Value = 4471367.5
$ws.Range("L99").Value = 1723.4286
$ws.Range("M99").Value = -4469869.5
$ws.Range("N99").Value = -4719.4286

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2981486
$ws.Range("I126").Value = 4471367.5
$ws.Range("J126").Value = 1723.4286
$ws.Range("K126").Value = 13414102.5
$ws.Range("L126").Value = 5170.2858
$ws.Range("M126").Value = -13411632.5
$ws.Range("N126").Value = -10110.2858

# CUL!row107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 592.9231
$ws.Range("J107").Value = 554.5
$ws.Range("L107").Value = 1663.5
$ws.Range("N107").Value = -5503.5

# CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 914.1429000000001
$ws.Range("J122").Value = 1159.8
$ws.Range("L122").Value = 10438.2
$ws.Range("N122").Value = -15338.2

# CUL!row125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 2827.5217
$ws.Range("I125").Value = 2500
$ws.Range("J125").Value = 2858.7144
$ws.Range("K125").Value = 7500
$ws.Range("L125").Value = 8576.143199999999
$ws.Range("M125").Value = -2580
$ws.Range("N125").Value = -18416.1432

# GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1353.6
$ws.Range("I113").Value = 962.6667
$ws.Range("J113").Value = 1940
$ws.Range("K113").Value = 962.6667
$ws.Range("L113").Value = 1940
$ws.Range("M113").Value = 1207.3333
$ws.Range("N113").Value = -6280

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 619631.5
$ws.Range("I122").Value = 1236356.8
$ws.Range("J122").Value = 2906.2222
$ws.Range("K122").Value = 3709070.4
$ws.Range("L122").Value = 8718.6666
$ws.Range("M122").Value = -3706620.4
$ws.Range("N122").Value = -13618.6666

# GSM!row126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2412
$ws.Range("I126").Value = 2426.6667
$ws.Range("J126").Value = 2409.6843
$ws.Range("K126").Value = 7280.000100000001
$ws.Range("L126").Value = 7229.0529
$ws.Range("M126").Value = -4810.000100000001
$ws.Range("N126").Value = -12169.0529

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3259.9714
$ws.Range("I132").Value = 3065.077
$ws.Range("J132").Value = 3823
$ws.Range("K132").Value = 9195.231
$ws.Range("L132").Value = 11469
$ws.Range("M132").Value = -6665.231
$ws.Range("N132").Value = -16529

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2661.0286
$ws.Range("I7").Value = 2012.4
$ws.Range("K7").Value = 2012.4
$ws.Range("M7").Value = -1900.4

# LTW!row16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4765697.5
$ws.Range("I16").Value = 9092241
$ws.Range("J16").Value = 6500
$ws.Range("K16").Value = 9092241
$ws.Range("L16").Value = 6500
$ws.Range("M16").Value = -9092071
$ws.Range("N16").Value = -6840

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 783.7143
$ws.Range("I22").Value = 840
$ws.Range("J22").Value = 752.44446
$ws.Range("K22").Value = 840
$ws.Range("L22").Value = 752.44446
$ws.Range("M22").Value = -545
$ws.Range("N22").Value = -1342.44446

# LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 783.7143
$ws.Range("I27").Value = 840
$ws.Range("J27").Value = 752.44446
$ws.Range("K27").Value = 840
$ws.Range("L27").Value = 752.44446
$ws.Range("M27").Value = -733
$ws.Range("N27").Value = -966.44446

# LTW!row40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2804.9524
$ws.Range("I40").Value = 1857.7142
$ws.Range("J40").Value = 3278.5715
$ws.Range("K40").Value = 1857.7142
$ws.Range("L40").Value = 3278.5715
$ws.Range("M40").Value = -1721.7142
$ws.Range("N40").Value = -3550.5715

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2661.0286
$ws.Range("I126").Value = 2012.4
$ws.Range("K126").Value = 6037.200000000001
$ws.Range("M126").Value = -3567.200000000001

# WVR!row16
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
